$d = $word.ActiveDocument

# --- Step 1: fix font for every paragraph first (cheap, no text edits yet) ---
$allParas = $d.Paragraphs
for ($i = 1; $i -le $allParas.Count; $i++) {
    try {
        $allParas.Item($i).Range.Font.Name = "Times New Roman"
    } catch {}
}

# --- Step 2: rewrite paragraph texts (this will coalesce runs within each paragraph) ---
# paragraph 1
$p = $d.Paragraphs.Item(1)
$pStart = $p.Range.Start
$pEnd = $p.Range.End - 1
$r = $d.Range($pStart, $pEnd)
$r.Text = "The Enchanting Realm of Biology: Unveiling the Wonders of Life"

# paragraph 2
$p = $d.Paragraphs.Item(2)
$pStart = $p.Range.Start
$pEnd = $p.Range.End - 1
$r = $d.Range($pStart, $pEnd)
$r.Text = "Professor Emily Carter"

# paragraph 3
$p = $d.Paragraphs.Item(3)
$pStart = $p.Range.Start
$pEnd = $p.Range.End - 1
$r = $d.Range($pStart, $pEnd)
$r.Text = "emilycarter@hscemail.com"

# paragraph 5
$p = $d.Paragraphs.Item(5)
$pStart = $p.Range.Start
$pEnd = $p.Range.End - 1
$r = $d.Range($pStart, $pEnd)
$r.Text = "Biology, the study of life, embarks on a thrilling voyage into the captivating realm of living organisms. It unveils the intricate intricacies of cells, uncovers the mysteries embedded within DNA, and delves into the extraordinary diversity of life forms inhabiting our planet. From the microscopic world of bacteria to the majestic grandeur of whales, biology unlocks the secrets of existence, revealing the profound interconnectedness that binds all living things. In this exploration of life's wonders, we embark on an exhilarating journey, unraveling the tapestry of life's rich symphony.`v`vThe study of genetics unveils the intricate mechanisms of heredity, deciphering the genetic blueprints that govern the traits and characteristics passed down through generations. We delve into the realm of evolution, tracing the grand narrative of life's transformation across eons, driven by the forces of natural selection. Through the lens of ecology, we unravel the delicate balance of ecosystems, revealing the intricate web of interactions that sustain the harmony of life.`v`vBiology, however, is not merely an intellectual pursuit; it holds profound implications for our daily lives and the future of our planet. It empowers us to understand the human body, enabling us to devise innovative treatments for diseases and ameliorate human suffering. By comprehending the intricate workings of ecosystems, we can devise strategies to protect and preserve the delicate balance of nature. As we continue to unravel the mysteries of life, biology empowers us to confront global challenges such as food security, climate change, and the preservation of biodiversity."

# paragraph 6
$p = $d.Paragraphs.Item(6)
$pStart = $p.Range.Start
$pEnd = $p.Range.End - 1
$r = $d.Range($pStart, $pEnd)
$r.Text = "Summary"

# paragraph 7
$p = $d.Paragraphs.Item(7)
$pStart = $p.Range.Start
$pEnd = $p.Range.End - 1
$r = $d.Range($pStart, $pEnd)
$r.Text = "Biology, the study of life, unveils the intricacies of living organisms, delving into the microscopic world of cells, deciphering the genetic blueprints of DNA, and exploring the magnificent diversity of life forms. By unraveling the tapestry of life's symphony, biology empowers us to understand the mechanisms of heredity, trace the narrative of evolution, and unravel the delicate balance of ecosystems. It holds profound implications for our daily lives and the future of our planet, enabling us to devise treatments for diseases, protect the environment, and confront global challenges. Biology's journey of discovery continues to inspire awe and wonder, revealing the interconnectedness of all living things and the profound beauty of the natural world."

# --- Step 3: re-split runs by re-applying font name at each desired run boundary ---
# paragraph 1 (single run)
$p = $d.Paragraphs.Item(1)
$pStart = $p.Range.Start
$len0 = "The Enchanting Realm of Biology: Unveiling the Wonders of Life".Length
$rr = $d.Range($pStart, $pStart + $len0)
$rr.Font.Name = "Times New Roman"

# paragraph 2 (single run)
$p = $d.Paragraphs.Item(2)
$pStart = $p.Range.Start
$len0 = "Professor Emily Carter".Length
$rr = $d.Range($pStart, $pStart + $len0)
$rr.Font.Name = "Times New Roman"

# paragraph 3 (3 runs)
$p = $d.Paragraphs.Item(3)
$pStart = $p.Range.Start
$off = 0
$seglen = "emilycarter@hscemail".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = "com".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen

# paragraph 5 (24 runs)
$p = $d.Paragraphs.Item(5)
$pStart = $p.Range.Start
$off = 0
$seglen = "Biology, the study of life, embarks on a thrilling voyage into the captivating realm of living organisms".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = " It unveils the intricate intricacies of cells, uncovers the mysteries embedded within DNA, and delves into the extraordinary diversity of life forms inhabiting our planet".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = " From the microscopic world of bacteria to the majestic grandeur of whales, biology unlocks the secrets of existence, revealing the profound interconnectedness that binds all living things".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = " In this exploration of life's wonders, we embark on an exhilarating journey, unraveling the tapestry of life's rich symphony".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = "`v".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = "`vThe study of genetics unveils the intricate mechanisms of heredity, deciphering the genetic blueprints that govern the traits and characteristics passed down through generations".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = " We delve into the realm of evolution, tracing the grand narrative of life's transformation across eons, driven by the forces of natural selection".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = " Through the lens of ecology, we unravel the delicate balance of ecosystems, revealing the intricate web of interactions that sustain the harmony of life".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = "`v".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = "`vBiology, however, is not merely an intellectual pursuit; it holds profound implications for our daily lives and the future of our planet".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = " It empowers us to understand the human body, enabling us to devise innovative treatments for diseases and ameliorate human suffering".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = " By comprehending the intricate workings of ecosystems, we can devise strategies to protect and preserve the delicate balance of nature".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = " As we continue to unravel the mysteries of life, biology empowers us to confront global challenges such as food security, climate change, and the preservation of biodiversity".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen

# paragraph 6 (single run)
$p = $d.Paragraphs.Item(6)
$pStart = $p.Range.Start
$len0 = "Summary".Length
$rr = $d.Range($pStart, $pStart + $len0)
$rr.Font.Name = "Times New Roman"

# paragraph 7 (8 runs)
$p = $d.Paragraphs.Item(7)
$pStart = $p.Range.Start
$off = 0
$seglen = "Biology, the study of life, unveils the intricacies of living organisms, delving into the microscopic world of cells, deciphering the genetic blueprints of DNA, and exploring the magnificent diversity of life forms".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = " By unraveling the tapestry of life's symphony, biology empowers us to understand the mechanisms of heredity, trace the narrative of evolution, and unravel the delicate balance of ecosystems".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = " It holds profound implications for our daily lives and the future of our planet, enabling us to devise treatments for diseases, protect the environment, and confront global challenges".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = " Biology's journey of discovery continues to inspire awe and wonder, revealing the interconnectedness of all living things and the profound beauty of the natural world".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen
$seglen = ".".Length
$rr = $d.Range($pStart + $off, $pStart + $off + $seglen)
$rr.Font.Name = "Times New Roman"
$off = $off + $seglen

